# Taches.xlsx - "Page d'accueil pour les veterinaires done" edit
#
# Row 13 (Voir les disponibilites des veterinaires / "Alejandro" / "Done")
# gets a Responsable and a Done? value, matching the same pattern already
# used on row 3 ("Alejandro" / "YES").
# The active selection is left on B9 (where the editor's cursor ended up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the Responsable (B13) and Done? (C13) cells for the new task row,
# reusing the exact same text already used elsewhere in the sheet so the
# shared-string table is reused rather than growing with new entries.
$ws.Range("B13").Value = "Alejandro"
$ws.Range("C13").Value = "YES"

# Leave the selection where the author ended up after the edit.
$ws.Range("B9").Select()
